# Adds a condensed "Post deployment" checklist table.
# 1) New table at J1:P14 (7 cols x 14 rows, header + 13 tasks) on the existing sheet.
# 2) The same table duplicated at A29:G42 (below the two existing tables).
# 3) Selection left on the new J1:P14 block (matches the authored workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels reused from the existing A1:G1 header (same six columns).
$headers = @("Task No", "Task Duration (mins)", "Start Date/Time", "End Date/Time", "Activity Description", "Activity Comments", "Task Dependencies")

$rows = @(
    @(1, 15, "Verify all PROD Databricks jobs and pipelines are running successfully", "No failed jobs", "200–250, 260–380", $false),
    @(2, 15, "Validate Event Hub ingestion to Bronze tables", "MQS events landing", "270–300", $false),
    @(3, 20, "Validate Silver & Gold table population", "Data flowing as expected", "310–380, 390–410", $false),
    @(4, 20, "Verify Data Reconciliation job execution", "Job completed successfully", "200, 210", $false),
    @(5, 15, "Validate reconciliation results against source", "No abnormal variance", "200, 210", $false),
    @(6, 15, "Verify Quote Cache Deletion pipeline execution", "Scheduled job triggered", 220, $true),
    @(7, 10, "Validate Quote Cache deletion effect", "Spot check quotes", 220, $true),
    @(8, 20, "Validate CoS ODS, Analyst & Pricing views", "Views accessible", "235, 240, 250", $false),
    @(9, 15, "Validate Vehicle & Area LRT service data", "LRT data available", "420–480", $false),
    @(10, 10, "Verify PV group permissions in PROD Databricks", "Access confirmed", "460, 490–510", $false),
    @(11, 15, "Monitor PROD logs and pipelines for errors", "No critical alerts", "200–510", $false),
    @(12, 10, "Business smoke testing on key tables/views", "High-level validation", "200–510", $false),
    @(13, 5, "Confirm monitoring and alerts status", "Stable system", "200–510", $false),
)

# --- Header row for the new J1:P1 block ---
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = 10 + $i   # J=10 .. P=16
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

# --- Data rows for the new J2:P14 block ---
foreach ($row in $rows) {
    $r = 1 + $row[0]   # task 1 -> row 2 ... task 13 -> row 14
    $ws.Cells.Item($r, 10).Value = $row[0]         # J: Task No
    $ws.Cells.Item($r, 11).Value = $row[1]         # K: Task Duration (mins)
    $ws.Cells.Item($r, 12).Value = "Post deployment"  # L: Start Date/Time
    $ws.Cells.Item($r, 13).Value = "Post deployment"  # M: End Date/Time
    $ws.Cells.Item($r, 14).Value = $row[2]         # N: Activity Description
    $ws.Cells.Item($r, 15).Value = $row[3]         # O: Activity Comments

    $depCell = $ws.Cells.Item($r, 16)              # P: Task Dependencies
    $depCell.Value = $row[4]
    $depCell.Font.Bold = $true
}

# --- Header row for the duplicated A29:G29 block ---
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = 1 + $i   # A=1 .. G=7
    $cell = $ws.Cells.Item(29, $col)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
}

# --- Data rows for the duplicated A30:G42 block ---
foreach ($row in $rows) {
    $r = 29 + $row[0]   # task 1 -> row 30 ... task 13 -> row 42
    $ws.Cells.Item($r, 1).Value = $row[0]          # A: Task No
    $ws.Cells.Item($r, 2).Value = $row[1]          # B: Task Duration (mins)
    $ws.Cells.Item($r, 3).Value = "Post deployment"   # C: Start Date/Time
    $ws.Cells.Item($r, 4).Value = "Post deployment"   # D: End Date/Time
    $ws.Cells.Item($r, 5).Value = $row[2]          # E: Activity Description
    $ws.Cells.Item($r, 6).Value = $row[3]          # F: Activity Comments

    $depCell2 = $ws.Cells.Item($r, 7)              # G: Task Dependencies
    $depCell2.Value = $row[4]
    $depCell2.Font.Bold = $true
}

# Match the authored selection: active cell J1, selection J1:P14.
$excel.Goto($ws.Range("J1:P14"))

